{"js": "// The document's narrative paragraphs (the free-text \"answer\" under each\n// Heading2/label) each got the text that used to belong to the *next*\n// content paragraph in this cyclical order:\n//   Objetivos -> Docente(s) -> Bibliografia -> Norma de recupera\u00e7\u00e3o ->\n//   Crit\u00e9rio -> M\u00e9todo -> Programa -> Programa resumido -> (back to Objetivos)\n// i.e. every paragraph's old text slides one step \"forward\" along that\n// cycle. We apply the 8 replacements directly against the known old/new\n// text values so the edit is self-contained and order-independent.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Simple paragraphs: the whole paragraph is a single run of free text, so\n// we can just overwrite the paragraph's full text.\nconst wholeParagraphEdits = [\n  {\n    oldText: \"Complementar a forma\u00e7\u00e3o multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, t\u00f3picos atuais e relevantes sobre gest\u00e3o da qualidade.\",\n    newText: \"A definir, de acordo com o t\u00f3pico programado.\",\n  },\n  {\n    oldText: \"5840535 - Messias Borges Silva\",\n    newText: \"Complementar a forma\u00e7\u00e3o multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, t\u00f3picos atuais e relevantes sobre gest\u00e3o da qualidade.\",\n  },\n  {\n    oldText: \"A definir, de acordo com o t\u00f3pico programado.\",\n    newText: \"O conte\u00fado desta disciplina ser\u00e1 de acordo com o t\u00f3pico a ser programado, devendo abordar assuntos complementares a forma\u00e7\u00e3o de um profissional de Engenharia.\",\n  },\n  {\n    oldText: \"O conte\u00fado desta disciplina ser\u00e1 de acordo com o t\u00f3pico a ser programado, devendo abordar assuntos complementares a forma\u00e7\u00e3o de um profissional de Engenharia.\",\n    newText: \"O desenvolvimento da disciplina ser\u00e1 baseado em leituras, aula expositiva, discuss\u00e3o e resolu\u00e7\u00e3o de estudos de caso e resolu\u00e7\u00e3o de exerc\u00edcios.\",\n  },\n  {\n    oldText: \"Textos fornecidos pelo professor da disciplina\\u000bArtigos extra\u00eddos de revistas especializadas na \u00e1rea de gest\u00e3o e produ\u00e7\u00e3o.\",\n    newText: \"5840535 - Messias Borges Silva\",\n  },\n];\n\n// Process paragraphs in document order, matching each against its current\n// (pre-edit) full text exactly once, so duplicate/overlapping literal\n// strings in the replacement chain can't cross-match each other.\nconst used = new Set();\nfor (const para of paragraphs.items) {\n  para.load(\"text\");\n}\nawait context.sync();\n\nfor (const para of paragraphs.items) {\n  const text = para.text.replace(/\\r/g, \"\");\n  for (let i = 0; i < wholeParagraphEdits.length; i++) {\n    if (used.has(i)) continue;\n    if (text === wholeParagraphEdits[i].oldText) {\n      para.insertText(wholeParagraphEdits[i].newText, \"Replace\");\n      used.add(i);\n      break;\n    }\n  }\n}\nawait context.sync();\n\n// Paragraphs that mix bold \"Label: \" runs with plain value runs: target the\n// specific value substring with search() so the bold labels and the\n// paragraph's other runs stay untouched. Each of these three values feeds\n// into the *next* one below it (M\u00e9todo -> Crit\u00e9rio -> Norma), so the edits\n// are applied from the bottom up: that way every search() still targets\n// pristine, not-yet-rewritten text, even though \"Provas e trabalhos.\" and\n// \"Prova \u00fanica...\" each exist as both an old value and a soon-to-be-written\n// new value within the same paragraph.\n//\n// Scope the search to the \"Avalia\u00e7\u00e3o\" paragraph itself (found by its bold\n// \"M\u00e9todo: \" label) instead of the whole body: by the time this runs,\n// step 1 has already copied some of these exact value strings into other,\n// unrelated paragraphs earlier in the document, and a body-wide search()\n// would latch onto the first (wrong) occurrence.\nconst avaliacaoPara = paragraphs.items.find((p) => p.text.replace(/\\r/g, \"\").indexOf(\"M\u00e9todo: \") === 0);\nif (!avaliacaoPara) {\n  throw new Error(\"Could not find the 'Avalia\u00e7\u00e3o' (M\u00e9todo/Crit\u00e9rio/Norma) paragraph\");\n}\n\nconst subRangeEdits = [\n  { oldText: \"Prova \u00fanica com nota maior ou igual a 5,0 (cinco).\", newText: \"Textos fornecidos pelo professor da disciplina\\u000bArtigos extra\u00eddos de revistas especializadas na \u00e1rea de gest\u00e3o e produ\u00e7\u00e3o.\" },\n  { oldText: \"Provas e trabalhos.\", newText: \"Prova \u00fanica com nota maior ou igual a 5,0 (cinco).\" },\n  { oldText: \"O desenvolvimento da disciplina ser\u00e1 baseado em leituras, aula expositiva, discuss\u00e3o e resolu\u00e7\u00e3o de estudos de caso e resolu\u00e7\u00e3o de exerc\u00edcios.\", newText: \"Provas e trabalhos.\" },\n];\n\nfor (const { oldText, newText } of subRangeEdits) {\n  const results = avaliacaoPara.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# The document's narrative paragraphs (the free-text \"answer\" under each\n# Heading2/label) each got the text that used to belong to the *next*\n# content paragraph in this cyclical order:\n#   Objetivos -> Docente(s) -> Bibliografia -> Norma de recupera\u00e7\u00e3o ->\n#   Crit\u00e9rio -> M\u00e9todo -> Programa -> Programa resumido -> (back to Objetivos)\n# i.e. every paragraph's old text slides one step \"forward\" along that\n# cycle. We apply the 8 replacements directly against the known old/new\n# text values so the edit is self-contained and order-independent.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: simple paragraphs where the whole paragraph is one run of\n# free text - overwrite the paragraph's full text in one shot.\n$wholeParagraphEdits = @(\n    @{\n        Old = \"Complementar a forma\u00e7\u00e3o multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, t\u00f3picos atuais e relevantes sobre gest\u00e3o da qualidade.\"\n        New = \"A definir, de acordo com o t\u00f3pico programado.\"\n    },\n    @{\n        Old = \"5840535 - Messias Borges Silva\"\n        New = \"Complementar a forma\u00e7\u00e3o multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, t\u00f3picos atuais e relevantes sobre gest\u00e3o da qualidade.\"\n    },\n    @{\n        Old = \"A definir, de acordo com o t\u00f3pico programado.\"\n        New = \"O conte\u00fado desta disciplina ser\u00e1 de acordo com o t\u00f3pico a ser programado, devendo abordar assuntos complementares a forma\u00e7\u00e3o de um profissional de Engenharia.\"\n    },\n    @{\n        Old = \"O conte\u00fado desta disciplina ser\u00e1 de acordo com o t\u00f3pico a ser programado, devendo abordar assuntos complementares a forma\u00e7\u00e3o de um profissional de Engenharia.\"\n        New = \"O desenvolvimento da disciplina ser\u00e1 baseado em leituras, aula expositiva, discuss\u00e3o e resolu\u00e7\u00e3o de estudos de caso e resolu\u00e7\u00e3o de exerc\u00edcios.\"\n    },\n    @{\n        Old = \"Textos fornecidos pelo professor da disciplina\" + [char]11 + \"Artigos extra\u00eddos de revistas especializadas na \u00e1rea de gest\u00e3o e produ\u00e7\u00e3o.\"\n        New = \"5840535 - Messias Borges Silva\"\n    }\n)\n\n# Snapshot every paragraph's current text up front (trimming the trailing\n# paragraph-mark \\r) so later writes in this same pass can't shadow one\n# another even though several old/new values repeat across the chain.\n$count = $d.Paragraphs.Count\n$snapshot = @()\nfor ($i = 1; $i -le $count; $i++) {\n    $raw = $d.Paragraphs.Item($i).Range.Text\n    $snapshot += $raw.TrimEnd([char]13)\n}\n\n$used = New-Object 'System.Collections.Generic.HashSet[int]'\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $snapshot[$i - 1]\n    for ($j = 0; $j -lt $wholeParagraphEdits.Count; $j++) {\n        if ($used.Contains($j)) { continue }\n        if ($text -ceq $wholeParagraphEdits[$j].Old) {\n            $d.Paragraphs.Item($i).Range.Text = $wholeParagraphEdits[$j].New\n            [void]$used.Add($j)\n            break\n        }\n    }\n}\n\n# --- Step 2: the \"Avalia\u00e7\u00e3o\" paragraph mixes bold \"Label: \" runs with\n# plain value runs, so only the value substrings are targeted (via a\n# Range.Find scoped to that paragraph) to leave the bold labels and w:br\n# breaks untouched. Each value feeds into the next one below it (M\u00e9todo ->\n# Crit\u00e9rio -> Norma), so we replace bottom-up: that way every Find still\n# matches pristine, not-yet-rewritten text, even though \"Provas e\n# trabalhos.\" and \"Prova \u00fanica...\" each exist as both an old value and a\n# soon-to-be-written new value within the same paragraph. Scoping to this\n# one paragraph (instead of searching $d.Content) also avoids matching the\n# identical strings Step 1 already copied into earlier, unrelated\n# paragraphs.\n$avaliacaoIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.StartsWith(\"M\u00e9todo: \")) {\n        $avaliacaoIndex = $i\n        break\n    }\n}\nif ($avaliacaoIndex -eq -1) {\n    throw \"Could not find the 'Avalia\u00e7\u00e3o' (M\u00e9todo/Crit\u00e9rio/Norma) paragraph\"\n}\n\n$subRangeEdits = @(\n    @{\n        Old = \"Prova \u00fanica com nota maior ou igual a 5,0 (cinco).\"\n        New = \"Textos fornecidos pelo professor da disciplina\" + [char]11 + \"Artigos extra\u00eddos de revistas especializadas na \u00e1rea de gest\u00e3o e produ\u00e7\u00e3o.\"\n    },\n    @{\n        Old = \"Provas e trabalhos.\"\n        New = \"Prova \u00fanica com nota maior ou igual a 5,0 (cinco).\"\n    },\n    @{\n        Old = \"O desenvolvimento da disciplina ser\u00e1 baseado em leituras, aula expositiva, discuss\u00e3o e resolu\u00e7\u00e3o de estudos de caso e resolu\u00e7\u00e3o de exerc\u00edcios.\"\n        New = \"Provas e trabalhos.\"\n    }\n)\n\nforeach ($edit in $subRangeEdits) {\n    $rng = $d.Paragraphs.Item($avaliacaoIndex).Range.Duplicate()\n    $found = $rng.Find.Execute($edit.Old, $false, $true, $false, $false, $false, $true, 1, $false, $edit.New, 1)\n    if (-not $found) {\n        throw \"Could not find text to replace: $($edit.Old)\"\n    }\n}\n"}
